# Applies the "cleaned up code, updated distributed.xlsx" revision:
#   - adds a small High/Low/Approx "demand" summary block (K6:Q8) plus a
#     header row (L5:Q5) to the "Normal capacity" sheet
#   - makes "Normal capacity" the active sheet again (it had drifted to
#     "Sheet1")

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Normal capacity")

# --- High demand row (taken straight from the '30-30 capacity' sheet) ----
$ws1.Range("K7").Value = "High demand"
$ws1.Range("L7").Value = 30
$ws1.Range("M7").Formula = "='30-30 capacity'!C2"
$ws1.Range("N7").Formula = "='30-30 capacity'!D2"
$ws1.Range("O7").Formula = "='30-30 capacity'!E2"
$ws1.Range("P7").Formula = "='30-30 capacity'!F2"
$ws1.Range("Q7").Formula = "='30-30 capacity'!G2"

# --- Low demand row (average over the Normal-capacity columns) ----------
$ws1.Range("K8").Value = "Low demand"
$ws1.Range("L8").Value = 20
$ws1.Range("M8").Formula = "=AVERAGE(C:C)"
$ws1.Range("N8").Formula = "=AVERAGE(D:D)"
$ws1.Range("O8").Formula = "=AVERAGE(E:E)"
$ws1.Range("P8").Formula = "=AVERAGE(F:F)"
$ws1.Range("Q8").Formula = "=AVERAGE(G:G)"

# --- High demand (Approx) row -- linear extrapolation of the two above --
$ws1.Range("K6").Value = "High demand (Approx)"
$ws1.Range("L6").Value = 40
$ws1.Range("M6").Formula = "=M7+M7-M8"
$ws1.Range("N6").Formula = "=N7+N7-N8"
$ws1.Range("O6").Formula = "=O7+O7-O8"
$ws1.Range("P6").Formula = "=P7+P7-P8"
$ws1.Range("Q6").Formula = "=Q7+Q7-Q8"
$ws1.Range("L6:Q6").Font.Italic = $true

# --- Header row for the little table (mirrors row 1's headers) ----------
$ws1.Range("L5").Value = "Numb. Of Aircraft"
$ws1.Range("M5").Formula = "=C1"
$ws1.Range("N5").Formula = "=D1"
$ws1.Range("O5").Formula = "=E1"
$ws1.Range("P5").Formula = "=F1"
$ws1.Range("Q5").Formula = "=G1"

# Size the new columns to their contents, like the existing C:H block.
$ws1.Columns("K:Q").AutoFit() | Out-Null

# Restore "Normal capacity" as the active sheet / selection.
$ws1.Activate() | Out-Null
$ws1.Range("N15").Select() | Out-Null
